$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.475.76'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '
$ws.Cells.Item(3, 4).Value = '1.878.09'
$ws.Cells.Item(3, 5).Value = '  +1.05%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.7162'
$ws.Cells.Item(5, 5).Value = '  +1.56%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '242.15'
$ws.Cells.Item(6, 5).Value = '  +1.50%  '
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07896'
$ws.Cells.Item(8, 5).Value = '  -1.49%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.3121'
$ws.Cells.Item(9, 5).Value = '  +2.94%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '25.24'
$ws.Cells.Item(10, 5).Value = '  +7.31%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08267'
$ws.Cells.Item(11, 5).Value = '  +0.92%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.7327'
$ws.Cells.Item(12, 5).Value = '  +3.60%  '
$ws.Cells.Item(13, 4).Value = '1.864.04'
$ws.Cells.Item(13, 5).Value = '  -1.25%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.285'
$ws.Cells.Item(14, 5).Value = '  +1.63%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '91.37'
$ws.Cells.Item(15, 5).Value = '  +1.80%  '
$ws.Cells.Item(16, 4).Value = '29.339.73'
$ws.Cells.Item(16, 5).Value = '  +0.29%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '5.931'
$ws.Cells.Item(17, 5).Value = '  +1.68%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '248.40'
$ws.Cells.Item(18, 5).Value = '  +4.17%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.000007900'
$ws.Cells.Item(19, 5).Value = '  -0.40%  '
$ws.Cells.Item(20, 5).Value = '  +0.16%  '
$ws.Cells.Item(21, 5).Value = '  +0.05%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.940'
$ws.Cells.Item(22, 5).Value = '  +6.21%  '
$ws.Cells.Item(23, 5).Value = '  +0.00%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.1590'
$ws.Cells.Item(24, 5).Value = '  +10.63%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '163.62'
$ws.Cells.Item(25, 5).Value = '  +0.44%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.034'
$ws.Cells.Item(26, 5).Value = '  +1.69%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.32'
$ws.Cells.Item(27, 5).Value = '  +1.09%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.362'
$ws.Cells.Item(28, 5).Value = '  -4.76%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.495'
$ws.Cells.Item(29, 5).Value = '  +1.35%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.378'
$ws.Cells.Item(30, 5).Value = '  +0.09%  '
$ws.Cells.Item(31, 5).Value = '  +2.50%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.05313'
$ws.Cells.Item(32, 5).Value = '  +2.07%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.936'
$ws.Cells.Item(33, 5).Value = '  +0.22%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.203'
$ws.Cells.Item(34, 5).Value = '  +3.55%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7257'
$ws.Cells.Item(35, 5).Value = '  +1.18%  '
$ws.Cells.Item(36, 5).Value = '  +0.42%  '
$ws.Cells.Item(37, 5).Value = '  +0.60%  '
$ws.Cells.Item(38, 4).Value = '1.267.82'
$ws.Cells.Item(38, 5).Value = '  +11.32%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.735'
$ws.Cells.Item(39, 5).Value = '  +0.17%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.9132'
$ws.Cells.Item(40, 5).Value = '  -2.68%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '74.37'
$ws.Cells.Item(41, 5).Value = '  +5.22%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '6.112'
$ws.Cells.Item(42, 5).Value = '  +2.13%  '
$ws.Cells.Item(43, 5).Value = '  +0.03%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '103.68'
$ws.Cells.Item(44, 5).Value = '  +0.76%  '
$ws.Cells.Item(45, 5).Value = '  +0.75%  '
$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.774'
$ws.Cells.Item(46, 5).Value = '  +0.51%  '
$ws.Cells.Item(47, 2).Value = 'SynthetixNetwork'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.938'
$ws.Cells.Item(47, 5).Value = '  +13.39%  '
$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.00000000120'
$ws.Cells.Item(48, 5).Value = '  -0.31%  '
$ws.Cells.Item(49, 2).Value = 'TheSandbox'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.4334'
$ws.Cells.Item(49, 5).Value = '  +1.48%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.282'
$ws.Cells.Item(50, 5).Value = '  +1.18%  '
$ws.Cells.Item(51, 2).Value = 'Aptos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.095'
$ws.Cells.Item(51, 5).Value = '  +1.65%  '
